{"js": "// Replace the inline \"Site Coverage Plan\" example picture with a plain\n// hyperlink run that points at the image's original URL on ura.gov.sg.\n//\n// Mirrors the authored diff: the <w:drawing> run inside the BodyText\n// paragraph is swapped for a <w:hyperlink> wrapping a run styled with the\n// built-in \"Hyperlink\" character style whose text is the external URL.\n\nconst body = context.document.body;\nconst pictures = body.inlinePictures;\npictures.load(\"items\");\nawait context.sync();\n\nif (pictures.items.length === 0) {\n  throw new Error(\"Expected an inline picture in the document body, found none.\");\n}\n\n// The document has exactly one inline image: the Site Coverage Plan example.\nconst picture = pictures.items[0];\nconst pictureRange = picture.getRange();\n\nconst url =\n  \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/SC01_Site_Coverage_Plan_Flats.jpg?h=100%25&w=100%25\";\n\n// Replace the picture (still inside its original BodyText paragraph) with\n// the URL text, then turn that new range into a hyperlink. Setting\n// `hyperlink` both creates the external relationship and applies the\n// built-in Hyperlink character style to the run, matching the target XML.\nconst linkRange = pictureRange.insertText(url, Word.InsertLocation.replace);\nlinkRange.hyperlink = url;\n\nawait context.sync();\n", "ps1": "# Replace the inline \"Site Coverage Plan\" example picture with a plain\n# hyperlink run that points at the image's original URL on ura.gov.sg.\n#\n# Mirrors the authored diff: the <w:drawing> run inside the BodyText\n# paragraph is swapped for a <w:hyperlink> wrapping a run styled with the\n# built-in \"Hyperlink\" character style whose text is the external URL.\n\n$d = $word.ActiveDocument\n$ils = $d.InlineShapes\n\nif ($ils.Count -eq 0) {\n    throw \"Expected an inline picture in the document body, found none.\"\n}\n\n# The document has exactly one inline image: the Site Coverage Plan example.\n$picture = $ils.Item(1)\n$picRange = $picture.Range\n\n$url = \"https://ura.gov.sg/-/media/Corporate/Guidelines/Development-control/Flats-Condominiums/SC01_Site_Coverage_Plan_Flats.jpg?h=100%25&w=100%25\"\n$start = $picRange.Start\n\n# Replace the picture (still inside its original BodyText paragraph) with\n# the URL text, then turn that new range into a hyperlink. Hyperlinks.Add\n# both creates the external relationship and applies the built-in Hyperlink\n# character style to the run, matching the target XML.\n$picRange.Text = $url\n$linkRange = $d.Range($start, $start + $url.Length)\n$d.Hyperlinks.Add($linkRange, $url) | Out-Null\n"}
